# "insert Data api completed"
#
# The mailing list on Sheet1 (A1:A6) drops two stale addresses:
#   - nikita@sourcesoftsolutions.com   (was A2)
#   - pandeyutkarsh2407@gmail.com      (was A5)
# leaving four rows (FilterMail header + 3 addresses). The remaining
# addresses shift up to fill the gaps, and the now-unused trailing rows
# are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "utkarshpandey2697@gmail.com"
$ws.Range("A3").Value = "arpit@sourcesoftsolutions.com"
$ws.Range("A4").Value = "utkarsh.sourcesoft@gmail.com"

$ws.Range("A5:A6").EntireRow.Delete()

$ws.Range("A4").Select()
